# Update Clock Tower Event Proposal pricing.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Event Budget")

# --- Unit-rate labels (text, column E) ---
# E8 and E9 originally share the same text ("400 / sqm"); set both so the
# rendered text stays identical for both line items, matching the source
# shared-string-table edit.
$ws.Range("E8").Value  = "450 / sqm"      # Center LED Screen
$ws.Range("E9").Value  = "450 / sqm"      # Side LED Screens
$ws.Range("E12").Value = "200 / sqm"      # Main Stage Structure
$ws.Range("E29").Value = "220 / unit"     # VIP Sofas
$ws.Range("E30").Value = "1,050 / table"  # Balcony Hospitality Tables
$ws.Range("E32").Value = "500 / person"   # Ushers

# --- Line item totals (column F) ---
$ws.Range("F8").Value  = 18000   # Center LED Screen
$ws.Range("F9").Value  = 10800   # Side LED Screens
$ws.Range("F10").Value = 6500    # Media Server
$ws.Range("F12").Value = 12800   # Main Stage Structure
$ws.Range("F13").Value = 4500    # Stage Carpet
$ws.Range("F15").Value = 16500   # Line Array System
$ws.Range("F21").Value = 8500    # Truss Structure
$ws.Range("F22").Value = 5500    # Stage & Ambient Lights
$ws.Range("F24").Value = 9000    # Guest Hospitality
$ws.Range("F26").Value = 10500   # Branded Photo Wall
$ws.Range("F27").Value = 4500    # Instant Photography
$ws.Range("F29").Value = 11000   # VIP Sofas
$ws.Range("F30").Value = 21000   # Balcony Hospitality Tables
$ws.Range("F32").Value = 5000    # Ushers

# --- Summary rows ---
$ws.Range("F36").Value = 146100  # SUBTOTAL
$ws.Range("F37").Value = 18900   # AGENCY COMMISSION
$ws.Range("F39").Value = 165000  # TOTAL PROJECT COST
